$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells for the ISO / VM path columns
$ws.Range("F1").Value = "ISO"
$ws.Range("G1").Value = "Path"

# Add new vSwitch column header
$ws.Range("I1").Value = "vSwitch"

# Populate the new vSwitch column for each VM data row
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = "External Switch"
}

# Last VM (Klient1) now belongs to site 3 and uses a different VM path
$ws.Range("A10").Value = 3
$ws.Range("G10").Value = "D:\Hyper-V"

# Resize the new vSwitch column to fit its contents
$ws.Columns.Item(9).ColumnWidth = 12.6

# Move the active selection as recorded in the workbook
$ws.Range("E15").Select()
